# Adapt column header formatting to respective input file names.
# "_old" header suffixes -> "_FV2210", "_new" header suffixes -> "_FV2304",
# then (re)build the Table1 listobject over the header row, and freeze the
# header row in the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (columns A1:J1 = "_old" -> "_FV2210", L1:U1 = "_new" -> "_FV2304") ---
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeadersFV2210 = @("Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210","Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210")

$oldHeaders2 = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")
$newHeadersFV2304 = @("Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2210[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2304[$i]
}

# --- 2. Freeze the header row (split after row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into a table so it matches the exported workbook ---
$rng = $ws.Range("A1:U63")
$tbl = $ws.ListObjects.Add(1, $rng, $false, 1, $null)
$tbl.Name = "Table1"
